# Scheduled runner: refresh market-price derived columns (H-N) across the
# Sheets workbook. Values come from the latest Universalis price pull.
$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# --- ALC ---------------------------------------------------------------
Set-Row "ALC" 16 @{
    H = 32625.8
    I = 3454.5
    J = 52073.332
    K = 3454.5
    L = 52073.332
    M = -3224.5
    N = -52533.332
}

Set-Row "ALC" 57 @{
    H = 20000
    J = 20000
    L = 60000
    N = -60998
}

Set-Row "ALC" 70 @{
    H = 2900.5
    I = 3100.2856
    J = 2745.111
    K = 9300.856800000001
    L = 8235.332999999999
    M = -9030.856800000001
    N = -8775.332999999999
}

Set-Row "ALC" 73 @{
    H = 2900.5
    I = 3100.2856
    J = 2745.111
    K = 9300.856800000001
    L = 8235.332999999999
    M = -8364.856800000001
    N = -10107.333
}

Set-Row "ALC" 101 @{
    H = 5235.4546
    J = 9173.333000000001
    L = 27519.999
    N = -30763.999
}

Set-Row "ALC" 106 @{
    H = 166670700
    I = 47623664
    K = 47623664
    M = -47623033
}

Set-Row "ALC" 116 @{
    H = 14425.375
    I = 34801
    K = 34801
    M = -31359
}

Set-Row "ALC" 118 @{
    H = 649.4
    I = 512.7143
    K = 1538.1429
    M = 118.8571000000002
}

Set-Row "ALC" 129 @{
    H = 955.0606
    I = 345.3846
    J = 1047.221
    K = 1036.1538
    L = 3141.663
    M = 3963.8462
    N = -13141.663
}

Set-Row "ALC" 138 @{
    H = 5334.48
    I = 1117
    J = 8146.1333
    K = 3351
    L = 24438.3999
    M = 1789
    N = -34718.3999
}

# --- ARM ---------------------------------------------------------------
Set-Row "ARM" 32 @{
    H = 6673.94
    I = 5609.8237
    K = 5609.8237
    M = -5322.8237
}

Set-Row "ARM" 61 @{
    H = 7846.3335
    I = 9493.643
    J = 2080.75
    K = 9493.643
    L = 2080.75
    M = -9281.643
    N = -2504.75
}

Set-Row "ARM" 74 @{
    H = 1673.8182
    I = 1591.6957
    J = 1862.7
    K = 1591.6957
    L = 1862.7
    M = -717.6957
    N = -3610.7
}

Set-Row "ARM" 77 @{
    H = 1673.8182
    I = 1591.6957
    J = 1862.7
    K = 7958.4785
    L = 9313.5
    M = -3590.4785
    N = -18049.5
}

Set-Row "ARM" 132 @{
    H = 6573.3335
    I = 1702.2727
    J = 9922.1875
    K = 5106.8181
    L = 29766.5625
    M = -2576.8181
    N = -34826.5625
}

Set-Row "ARM" 136 @{
    H = 7846.3335
    I = 9493.643
    J = 2080.75
    K = 28480.929
    L = 6242.25
    M = -25930.929
    N = -11342.25
}

# --- BSM ---------------------------------------------------------------
Set-Row "BSM" 64 @{
    H = 233.27272
    I = 181.5
    K = 181.5
    M = 43.5
}

Set-Row "BSM" 67 @{
    H = 233.27272
    I = 181.5
    K = 181.5
    M = 598.5
}

Set-Row "BSM" 107 @{
    H = 888.3333
    J = 824.55554
    L = 824.55554
    N = -4664.55554
}

Set-Row "BSM" 134 @{
    H = 5271
    I = 6352.696
    J = 3009.2727
    K = 19058.088
    L = 9027.8181
    M = -16523.088
    N = -14097.8181
}

# --- CRP ---------------------------------------------------------------
Set-Row "CRP" 22 @{
    H = 367.14285
    I = 328.33334
    J = 600
    K = 328.33334
    L = 600
    M = 21.66665999999998
    N = -1300
}

# --- GSM ---------------------------------------------------------------
Set-Row "GSM" 7 @{
    H = 5405715
    I = 9250000
    J = 280002
    K = 9250000
    L = 280002
    M = -9249888
    N = -280226
}

Set-Row "GSM" 8 @{
    H = 5405715
    I = 9250000
    J = 280002
    K = 9250000
    L = 280002
    M = -9249861
    N = -280280
}

Set-Row "GSM" 132 @{
    H = 3046.2222
    I = 2070.2307
    J = 3597.8696
    K = 6210.6921
    L = 10793.6088
    M = -3680.6921
    N = -15853.6088
}

# --- LTW ---------------------------------------------------------------
Set-Row "LTW" 20 @{
    H = 3000
    J = 3000
    L = 3000
    N = -3452
}

Set-Row "LTW" 21 @{
    H = 70007
    J = 70007
    L = 70007
    N = -70355
}

Set-Row "LTW" 24 @{
    H = 19633.334
    J = 19633.334
    L = 19633.334
    N = -20319.334
}

# --- WVR ---------------------------------------------------------------
Set-Row "WVR" 14 @{
    H = 53253.75
    I = 3000
    K = 3000
    M = -2832
}

Set-Row "WVR" 137 @{
    H = 59800
    J = 59800
    L = 59800
    N = -70000
}
